# Update LR-pair stats with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs target)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08856766666666667
$ws.Range("H2").Value = 0.265703
$ws.Range("M2").Value = 34.49888633333333
$ws.Range("N2").Value = 103.496659
$ws.Range("O2").Value = 0.4998067520528027
$ws.Range("P2").Value = 0.4998067520528027
$ws.Range("Q2").Value = 3.055485865141889
$ws.Range("R2").Value = 27.499372786277
$ws.Range("S2").Value = 0.4998067520528027
$ws.Range("T2").Value = 0.4998067520528027

# Row 3 (FAPs target)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08856766666666667
$ws.Range("H3").Value = 0.265703
$ws.Range("N3").Value = 72.35583600000001
$ws.Range("O3").Value = 0.3494212830891987
$ws.Range("P3").Value = 0.3494212830891987
$ws.Range("Q3").Value = 2.136129188078667
$ws.Range("R3").Value = 19.225162692708
$ws.Range("S3").Value = 0.3494212830891987
$ws.Range("T3").Value = 0.3494212830891987

# Row 4 (MuSCs target)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08856766666666667
$ws.Range("H4").Value = 0.265703
$ws.Range("M4").Value = 10.406952
$ws.Range("N4").Value = 31.220856
$ws.Range("O4").Value = 0.1507719648579985
$ws.Range("P4").Value = 0.1507719648579985
$ws.Range("Q4").Value = 0.9217194557519999
$ws.Range("R4").Value = 8.295475101768
$ws.Range("S4").Value = 0.1507719648579985
$ws.Range("T4").Value = 0.1507719648579985
